$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.043216674681025
$ws.Range("D2").Value = 1.051298328604505
$ws.Range("E2").Value = 1.057597504576344
$ws.Range("F2").Value = 1.064579297132833
$ws.Range("I2").Value = 1.047099349186998
$ws.Range("J2").Value = 1.048288311811995
$ws.Range("K2").Value = 1.054050218260536
$ws.Range("L2").Value = 1.06033203795109
$ws.Range("M2").Value = 1.067294853150819
$ws.Range("N2").Value = 1.049777001455118

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.044046547712729
$ws.Range("D3").Value = 1.051954060995268
$ws.Range("E3").Value = 1.058452524966801
$ws.Range("F3").Value = 1.065398262614001
$ws.Range("I3").Value = 1.047329446234386
$ws.Range("J3").Value = 1.048765535286832
$ws.Range("K3").Value = 1.054518972720222
$ws.Range("L3").Value = 1.061000821452982
$ws.Range("M3").Value = 1.067929041701335
$ws.Range("N3").Value = 1.050254902642027

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.044584025395771
$ws.Range("D4").Value = 1.052378801606131
$ws.Range("E4").Value = 1.059007083113257
$ws.Range("F4").Value = 1.065929136597933
$ws.Range("I4").Value = 1.047477299817973
$ws.Range("J4").Value = 1.049074143711101
$ws.Range("K4").Value = 1.054822039635879
$ws.Range("L4").Value = 1.06143422455105
$ws.Range("M4").Value = 1.068339697535765
$ws.Range("N4").Value = 1.050563949325657

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.044810097734625
$ws.Range("D5").Value = 1.052557465640658
$ws.Range("E5").Value = 1.059240529326697
$ws.Range("F5").Value = 1.066152541183709
$ws.Range("I5").Value = 1.047539208963685
$ws.Range("J5").Value = 1.049203836665245
$ws.Range("K5").Value = 1.054949388164897
$ws.Range("L5").Value = 1.061616582534054
$ws.Range("M5").Value = 1.068512405909791
$ws.Range("N5").Value = 1.050693826458672

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.044848063072463
$ws.Range("D6").Value = 1.052587470111943
$ws.Range("E6").Value = 1.059279744094052
$ws.Range("F6").Value = 1.066190064944435
$ws.Range("I6").Value = 1.047549589185722
$ws.Range("J6").Value = 1.049225609947768
$ws.Range("K6").Value = 1.054970766940137
$ws.Range("L6").Value = 1.061647210303491
$ws.Range("M6").Value = 1.068541408385911
$ws.Range("N6").Value = 1.050715630661754

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.044587045728826
$ws.Range("D7").Value = 1.05238118852087
$ws.Range("E7").Value = 1.059010201217433
$ws.Range("F7").Value = 1.065932120858622
$ws.Range("I7").Value = 1.047478128028466
$ws.Range("J7").Value = 1.049075876857114
$ws.Range("K7").Value = 1.054823741513694
$ws.Range("L7").Value = 1.061436660619045
$ws.Range("M7").Value = 1.068342005005714
$ws.Range("N7").Value = 1.050565684932936

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.043497030823613
$ws.Range("D8").Value = 1.051519844710667
$ws.Range("E8").Value = 1.057886191817061
$ws.Range("F8").Value = 1.064855872790885
$ws.Range("I8").Value = 1.047177325288805
$ws.Range("J8").Value = 1.048449629708434
$ws.Range("K8").Value = 1.05420868671295
$ws.Range("L8").Value = 1.060557919791478
$ws.Range("M8").Value = 1.067509118428362
$ws.Range("N8").Value = 1.049938548441474

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.041580143241208
$ws.Range("D9").Value = 1.050005478204498
$ws.Range("E9").Value = 1.055915606304584
$ws.Range("F9").Value = 1.062966736027776
$ws.Range("I9").Value = 1.046639384864262
$ws.Range("J9").Value = 1.047344724324417
$ws.Range("K9").Value = 1.053123035133944
$ws.Range("L9").Value = 1.059014556758687
$ws.Range("M9").Value = 1.066043783446044
$ws.Range("N9").Value = 1.048832073965076

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.04030490869416
$ws.Range("D10").Value = 1.048998312102363
$ws.Range("E10").Value = 1.054608762817408
$ws.Range("F10").Value = 1.061712367859601
$ws.Range("I10").Value = 1.046275503467192
$ws.Range("J10").Value = 1.046607268208597
$ws.Range("K10").Value = 1.052398103139715
$ws.Range("L10").Value = 1.057989165943723
$ws.Range("M10").Value = 1.065068546354015
$ws.Range("N10").Value = 1.048093570576985

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.039753376494002
$ws.Range("D11").Value = 1.048562792121035
$ws.Range("E11").Value = 1.054044541276425
$ws.Range("F11").Value = 1.061170436116738
$ws.Range("I11").Value = 1.046116702753811
$ws.Range("J11").Value = 1.046287754910638
$ws.Range("K11").Value = 1.05208393968676
$ws.Range("L11").Value = 1.057546014910915
$ws.Range("M11").Value = 1.064646669946485
$ws.Range("N11").Value = 1.047773603533494

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.039548612571521
$ws.Range("D12").Value = 1.048401110971235
$ws.Range("E12").Value = 1.053835214183206
$ws.Range("F12").Value = 1.060969323483392
$ws.Range("I12").Value = 1.046057531958352
$ws.Range("J12").Value = 1.046169046154073
$ws.Range("K12").Value = 1.051967207164077
$ws.Range("L12").Value = 1.057381538210117
$ws.Range("M12").Value = 1.064490029125157
$ws.Range("N12").Value = 1.047654726196879

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.039592530590569
$ws.Range("D13").Value = 1.048435788008979
$ws.Range("E13").Value = 1.053880104214768
$ws.Range("F13").Value = 1.061012454425271
$ws.Range("I13").Value = 1.046070232659857
$ws.Range("J13").Value = 1.046194510800339
$ws.Range("K13").Value = 1.05199224839545
$ws.Range("L13").Value = 1.057416813142446
$ws.Range("M13").Value = 1.064523626242646
$ws.Range("N13").Value = 1.047680227005864

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.039736448601184
$ws.Range("D14").Value = 1.048549425653073
$ws.Range("E14").Value = 1.054027233124149
$ws.Range("F14").Value = 1.061153808299552
$ws.Range("I14").Value = 1.046111815449531
$ws.Range("J14").Value = 1.046277942957329
$ws.Range("K14").Value = 1.052074291312332
$ws.Range("L14").Value = 1.057532416567933
$ws.Range("M14").Value = 1.064633720675685
$ws.Range("N14").Value = 1.047763777646086

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.039825134496349
$ws.Range("D15").Value = 1.048619453583385
$ws.Range("E15").Value = 1.054117917275329
$ws.Range("F15").Value = 1.061240925652509
$ws.Range("I15").Value = 1.04613741146995
$ws.Range("J15").Value = 1.046329344683702
$ws.Range("K15").Value = 1.052124835634529
$ws.Range("L15").Value = 1.057603660839154
$ws.Range("M15").Value = 1.064701561867232
$ws.Range("N15").Value = 1.047815252368806

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.040341525961419
$ws.Range("D16").Value = 1.049027228703426
$ws.Range("E16").Value = 1.054646243281069
$ws.Range("F16").Value = 1.061748359912894
$ws.Range("I16").Value = 1.046286016540041
$ws.Range("J16").Value = 1.04662846935839
$ws.Range("K16").Value = 1.05241894770101
$ws.Range("L16").Value = 1.058018594464837
$ws.Range("M16").Value = 1.065096553628715
$ws.Range("N16").Value = 1.048114801834843

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.040665620594684
$ws.Range("D17").Value = 1.049283174372239
$ws.Range("E17").Value = 1.054978091451378
$ws.Range("F17").Value = 1.062066987483637
$ws.Range("I17").Value = 1.04637890178391
$ws.Range("J17").Value = 1.046816052391576
$ws.Range("K17").Value = 1.052603366905601
$ws.Range("L17").Value = 1.058279100072391
$ws.Range("M17").Value = 1.065344431811463
$ws.Range("N17").Value = 1.048302651257451

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.040854722385345
$ws.Range("D18").Value = 1.049432519917693
$ws.Range("E18").Value = 1.055171811973963
$ws.Range("F18").Value = 1.062252954938014
$ws.Range("I18").Value = 1.046432960704105
$ws.Range("J18").Value = 1.046925447991928
$ws.Range("K18").Value = 1.052710910033354
$ws.Range("L18").Value = 1.058431130485196
$ws.Range("M18").Value = 1.065489054189796
$ws.Range("N18").Value = 1.048412202212104

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.040919211824409
$ws.Range("D19").Value = 1.049483452463636
$ws.Range("E19").Value = 1.055237892592092
$ws.Range("F19").Value = 1.062316384884641
$ws.Range("I19").Value = 1.046451373110037
$ws.Range("J19").Value = 1.046962745904243
$ws.Range("K19").Value = 1.052747575100283
$ws.Range("L19").Value = 1.058482982769173
$ws.Range("M19").Value = 1.065538373288223
$ws.Range("N19").Value = 1.048449553091732

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.040630841792641
$ws.Range("D20").Value = 1.049255707939659
$ws.Range("E20").Value = 1.054942470812701
$ws.Range("F20").Value = 1.062032789605115
$ws.Range("I20").Value = 1.046368948427939
$ws.Range("J20").Value = 1.046795928402558
$ws.Range("K20").Value = 1.052583583095789
$ws.Range("L20").Value = 1.058251141803126
$ws.Range("M20").Value = 1.065317832759186
$ws.Range("N20").Value = 1.048282498690061

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.03969406556905
$ws.Range("D21").Value = 1.04851595969046
$ws.Range("E21").Value = 1.053983900412909
$ws.Range("F21").Value = 1.061112177985575
$ws.Range("I21").Value = 1.046099575458685
$ws.Range("J21").Value = 1.046253375009971
$ws.Range("K21").Value = 1.052050132763003
$ws.Range("L21").Value = 1.057498370657416
$ws.Range("M21").Value = 1.064601298861001
$ws.Range("N21").Value = 1.047739174809425

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.039105654778151
$ws.Range("D22").Value = 1.048051374298421
$ws.Range("E22").Value = 1.053382656027918
$ws.Range("F22").Value = 1.060534423929619
$ws.Range("I22").Value = 1.045929139219208
$ws.Range("J22").Value = 1.045912092874709
$ws.Range("K22").Value = 1.05171451083456
$ws.Range("L22").Value = 1.057025822595749
$ws.Range("M22").Value = 1.064151150048306
$ws.Range("N22").Value = 1.047397408014386

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.039417527179314
$ws.Range("D23").Value = 1.048297609569836
$ws.Range("E23").Value = 1.053701249248513
$ws.Range("F23").Value = 1.060840600111948
$ws.Range("I23").Value = 1.046019591920943
$ws.Range("J23").Value = 1.046093027485078
$ws.Range("K23").Value = 1.051892450838092
$ws.Range("L23").Value = 1.057276257720748
$ws.Range("M23").Value = 1.06438974736418
$ws.Range("N23").Value = 1.047578599572654

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.040646556642168
$ws.Range("D24").Value = 1.049268118663089
$ws.Range("E24").Value = 1.05495856575471
$ws.Range("F24").Value = 1.062048241791944
$ws.Range("I24").Value = 1.046373446290304
$ws.Range("J24").Value = 1.046805021624001
$ws.Range("K24").Value = 1.05259252262673
$ws.Range("L24").Value = 1.058263774688708
$ws.Range("M24").Value = 1.06532985160496
$ws.Range("N24").Value = 1.048291604824921

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.042075237493477
$ws.Range("D25").Value = 1.050396560835878
$ws.Range("E25").Value = 1.05642384560634
$ws.Range("F25").Value = 1.063454240346802
$ws.Range("I25").Value = 1.046779385062849
$ws.Range("J25").Value = 1.047630524525224
$ws.Range("K25").Value = 1.053403912951063
$ws.Range("L25").Value = 1.059412940022482
$ws.Range("M25").Value = 1.066422323261479
$ws.Range("N25").Value = 1.049118280034951

Write-Host "vm_pu values updated for 380 kV case"
